$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(43, 8).Value = 5997.5
$ws.Cells.Item(43, 9).Value = 0
$ws.Cells.Item(43, 10).Value = 5997.5
$ws.Cells.Item(43, 11).Value = 0
$ws.Cells.Item(43, 12).Value = 5997.5
$ws.Cells.Item(43, 13).ClearContents()
$ws.Cells.Item(43, 14).Value = -6135.5

$ws.Cells.Item(51, 8).Value = 11928.286
$ws.Cells.Item(51, 9).Value = 13499.333
$ws.Cells.Item(51, 10).Value = 10750
$ws.Cells.Item(51, 11).Value = 13499.333
$ws.Cells.Item(51, 12).Value = 10750
$ws.Cells.Item(51, 13).Value = -13015.333
$ws.Cells.Item(51, 14).Value = -11718

$ws.Cells.Item(74, 8).Value = 5162
$ws.Cells.Item(74, 9).Value = 5162
$ws.Cells.Item(74, 11).Value = 5162
$ws.Cells.Item(74, 13).Value = -4226

$ws.Cells.Item(77, 8).Value = 5162
$ws.Cells.Item(77, 9).Value = 5162
$ws.Cells.Item(77, 11).Value = 25810
$ws.Cells.Item(77, 13).Value = -21130

$ws.Cells.Item(80, 8).Value = 3010
$ws.Cells.Item(80, 9).Value = 0
$ws.Cells.Item(80, 10).Value = 3010
$ws.Cells.Item(80, 11).Value = 0
$ws.Cells.Item(80, 12).Value = 9030
$ws.Cells.Item(80, 13).ClearContents()
$ws.Cells.Item(80, 14).Value = -11026

$ws.Cells.Item(83, 8).Value = 3010
$ws.Cells.Item(83, 9).Value = 0
$ws.Cells.Item(83, 10).Value = 3010
$ws.Cells.Item(83, 11).Value = 0
$ws.Cells.Item(83, 12).Value = 27090
$ws.Cells.Item(83, 13).ClearContents()
$ws.Cells.Item(83, 14).Value = -37074

$ws.Cells.Item(94, 8).Value = 3999
$ws.Cells.Item(94, 9).Value = 0
$ws.Cells.Item(94, 10).Value = 3999
$ws.Cells.Item(94, 11).Value = 0
$ws.Cells.Item(94, 12).Value = 3999
$ws.Cells.Item(94, 13).ClearContents()
$ws.Cells.Item(94, 14).Value = -4901

$ws.Cells.Item(111, 8).Value = 5299.6665
$ws.Cells.Item(111, 9).Value = 5299.6665
$ws.Cells.Item(111, 11).Value = 15898.9995
$ws.Cells.Item(111, 13).Value = -12831.9995

$ws.Cells.Item(138, 8).Value = 1898.5714
$ws.Cells.Item(138, 9).Value = 1706.1538
$ws.Cells.Item(138, 11).Value = 5118.4614
$ws.Cells.Item(138, 13).Value = 21.53859999999986

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(15, 8).Value = 12999.857
$ws.Cells.Item(15, 9).Value = 999
$ws.Cells.Item(15, 11).Value = 999
$ws.Cells.Item(15, 13).Value = -649

$ws.Cells.Item(39, 8).Value = 35000
$ws.Cells.Item(39, 9).Value = 35000
$ws.Cells.Item(39, 11).Value = 35000
$ws.Cells.Item(39, 13).Value = -34480

$ws.Cells.Item(61, 8).Value = 5068.1816
$ws.Cells.Item(61, 9).Value = 4439.25
$ws.Cells.Item(61, 11).Value = 4439.25
$ws.Cells.Item(61, 13).Value = -4227.25

$ws.Cells.Item(74, 8).Value = 699.6667
$ws.Cells.Item(74, 9).Value = 699.6667
$ws.Cells.Item(74, 10).Value = 0
$ws.Cells.Item(74, 11).Value = 699.6667
$ws.Cells.Item(74, 12).Value = 0
$ws.Cells.Item(74, 13).Value = 174.3333
$ws.Cells.Item(74, 14).ClearContents()

$ws.Cells.Item(77, 8).Value = 699.6667
$ws.Cells.Item(77, 9).Value = 699.6667
$ws.Cells.Item(77, 10).Value = 0
$ws.Cells.Item(77, 11).Value = 3498.3335
$ws.Cells.Item(77, 12).Value = 0
$ws.Cells.Item(77, 13).Value = 869.6665000000003
$ws.Cells.Item(77, 14).ClearContents()

$ws.Cells.Item(97, 8).Value = 1641.8182
$ws.Cells.Item(97, 9).Value = 826.8421
$ws.Cells.Item(97, 11).Value = 826.8421
$ws.Cells.Item(97, 13).Value = -330.8421

$ws.Cells.Item(136, 8).Value = 5068.1816
$ws.Cells.Item(136, 9).Value = 4439.25
$ws.Cells.Item(136, 11).Value = 13317.75
$ws.Cells.Item(136, 13).Value = -10767.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 4102.294
$ws.Cells.Item(86, 9).Value = 2955.7
$ws.Cells.Item(86, 11).Value = 2955.7
$ws.Cells.Item(86, 13).Value = -1832.7

$ws.Cells.Item(89, 8).Value = 4102.294
$ws.Cells.Item(89, 9).Value = 2955.7
$ws.Cells.Item(89, 11).Value = 14778.5
$ws.Cells.Item(89, 13).Value = -9162.5

$ws.Cells.Item(94, 8).Value = 910.1786
$ws.Cells.Item(94, 9).Value = 855.3570999999999
$ws.Cells.Item(94, 11).Value = 855.3570999999999
$ws.Cells.Item(94, 13).Value = -404.3570999999999

$ws.Cells.Item(99, 8).Value = 2178
$ws.Cells.Item(99, 9).Value = 1972.5
$ws.Cells.Item(99, 10).Value = 3000
$ws.Cells.Item(99, 11).Value = 1972.5
$ws.Cells.Item(99, 12).Value = 3000
$ws.Cells.Item(99, 13).Value = -474.5
$ws.Cells.Item(99, 14).Value = -5996

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 592.3125
$ws.Cells.Item(7, 9).Value = 499.7
$ws.Cells.Item(7, 10).Value = 746.6667
$ws.Cells.Item(7, 11).Value = 499.7
$ws.Cells.Item(7, 12).Value = 746.6667
$ws.Cells.Item(7, 13).Value = -386.7
$ws.Cells.Item(7, 14).Value = -972.6667

$ws.Cells.Item(16, 8).Value = 3078
$ws.Cells.Item(16, 9).Value = 766.3333
$ws.Cells.Item(16, 11).Value = 766.3333
$ws.Cells.Item(16, 13).Value = -479.3333

$ws.Cells.Item(62, 8).Value = 2282.5
$ws.Cells.Item(62, 10).Value = 2499.6667
$ws.Cells.Item(62, 12).Value = 2499.6667
$ws.Cells.Item(62, 14).Value = -3747.6667

$ws.Cells.Item(65, 8).Value = 2282.5
$ws.Cells.Item(65, 10).Value = 2499.6667
$ws.Cells.Item(65, 12).Value = 12498.3335
$ws.Cells.Item(65, 14).Value = -18738.3335

$ws.Cells.Item(99, 8).Value = 4833.3184
$ws.Cells.Item(99, 9).Value = 4916.65
$ws.Cells.Item(99, 11).Value = 4916.65
$ws.Cells.Item(99, 13).Value = -3418.65

$ws.Cells.Item(113, 8).Value = 3078
$ws.Cells.Item(113, 9).Value = 766.3333
$ws.Cells.Item(113, 11).Value = 766.3333
$ws.Cells.Item(113, 13).Value = 1403.6667

$ws.Cells.Item(126, 8).Value = 4833.3184
$ws.Cells.Item(126, 9).Value = 4916.65
$ws.Cells.Item(126, 11).Value = 14749.95
$ws.Cells.Item(126, 13).Value = -12279.95

$ws.Cells.Item(141, 8).Value = 62500
$ws.Cells.Item(141, 10).Value = 0
$ws.Cells.Item(141, 12).Value = 0
$ws.Cells.Item(141, 14).ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 38.133335
$ws.Cells.Item(12, 9).Value = 3
$ws.Cells.Item(12, 10).Value = 55.7
$ws.Cells.Item(12, 11).Value = 9
$ws.Cells.Item(12, 12).Value = 167.1
$ws.Cells.Item(12, 13).Value = 164
$ws.Cells.Item(12, 14).Value = -513.1

$ws.Cells.Item(23, 8).Value = 276.41666
$ws.Cells.Item(23, 9).Value = 275.66666
$ws.Cells.Item(23, 10).Value = 277.16666
$ws.Cells.Item(23, 11).Value = 826.9999799999999
$ws.Cells.Item(23, 12).Value = 831.4999799999999
$ws.Cells.Item(23, 13).Value = -591.9999799999999
$ws.Cells.Item(23, 14).Value = -1301.49998

$ws.Cells.Item(26, 8).Value = 581
$ws.Cells.Item(26, 9).Value = 413.75
$ws.Cells.Item(26, 11).Value = 1241.25
$ws.Cells.Item(26, 13).Value = -953.25

$ws.Cells.Item(138, 8).Value = 6930.3
$ws.Cells.Item(138, 9).Value = 6930.3
$ws.Cells.Item(138, 11).Value = 20790.9
$ws.Cells.Item(138, 13).Value = -15650.9

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(9, 8).Value = 621.55554
$ws.Cells.Item(9, 10).Value = 950
$ws.Cells.Item(9, 12).Value = 950
$ws.Cells.Item(9, 14).Value = -1290

$ws.Cells.Item(14, 8).Value = 11464379
$ws.Cells.Item(14, 10).Value = 175
$ws.Cells.Item(14, 12).Value = 175
$ws.Cells.Item(14, 14).Value = -511

$ws.Cells.Item(113, 8).Value = 3809.75
$ws.Cells.Item(113, 9).Value = 3809.75
$ws.Cells.Item(113, 11).Value = 3809.75
$ws.Cells.Item(113, 13).Value = -1639.75

$ws.Cells.Item(132, 8).Value = 3730
$ws.Cells.Item(132, 9).Value = 3684
$ws.Cells.Item(132, 11).Value = 11052
$ws.Cells.Item(132, 13).Value = -8522

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(14, 8).Value = 180
$ws.Cells.Item(14, 9).Value = 180
$ws.Cells.Item(14, 10).Value = 0
$ws.Cells.Item(14, 11).Value = 180
$ws.Cells.Item(14, 12).Value = 0
$ws.Cells.Item(14, 13).Value = -8
$ws.Cells.Item(14, 14).ClearContents()

$ws.Cells.Item(30, 8).Value = 907.75
$ws.Cells.Item(30, 9).Value = 907.75
$ws.Cells.Item(30, 11).Value = 907.75
$ws.Cells.Item(30, 13).Value = -799.75

$ws.Cells.Item(46, 8).Value = 1964.8
$ws.Cells.Item(46, 9).Value = 1324.125
$ws.Cells.Item(46, 11).Value = 1324.125
$ws.Cells.Item(46, 13).Value = -1136.125

$ws.Cells.Item(68, 8).Value = 3499.5
$ws.Cells.Item(68, 9).Value = 2999
$ws.Cells.Item(68, 11).Value = 2999
$ws.Cells.Item(68, 13).Value = -2250

$ws.Cells.Item(71, 8).Value = 3499.5
$ws.Cells.Item(71, 9).Value = 2999
$ws.Cells.Item(71, 11).Value = 14995
$ws.Cells.Item(71, 13).Value = -11251

$ws.Cells.Item(82, 8).Value = 2958.5
$ws.Cells.Item(82, 9).Value = 0
$ws.Cells.Item(82, 10).Value = 2958.5
$ws.Cells.Item(82, 11).Value = 0
$ws.Cells.Item(82, 12).Value = 2958.5
$ws.Cells.Item(82, 13).ClearContents()
$ws.Cells.Item(82, 14).Value = -3680.5

$ws.Cells.Item(85, 8).Value = 2958.5
$ws.Cells.Item(85, 9).Value = 0
$ws.Cells.Item(85, 10).Value = 2958.5
$ws.Cells.Item(85, 11).Value = 0
$ws.Cells.Item(85, 12).Value = 2958.5
$ws.Cells.Item(85, 13).ClearContents()
$ws.Cells.Item(85, 14).Value = -5454.5

$ws.Cells.Item(100, 8).Value = 1003
$ws.Cells.Item(100, 9).Value = 1003
$ws.Cells.Item(100, 11).Value = 1003
$ws.Cells.Item(100, 13).Value = -462

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 5912.8887
$ws.Cells.Item(96, 9).Value = 5888.3335
$ws.Cells.Item(96, 11).Value = 5888.3335
$ws.Cells.Item(96, 13).Value = -4515.3335

$ws.Cells.Item(122, 8).Value = 2846.7144
$ws.Cells.Item(122, 9).Value = 2813.818
$ws.Cells.Item(122, 10).Value = 2967.3333
$ws.Cells.Item(122, 11).Value = 8441.454000000002
$ws.Cells.Item(122, 12).Value = 8901.999899999999
$ws.Cells.Item(122, 13).Value = -5991.454000000002
$ws.Cells.Item(122, 14).Value = -13801.9999

$ws.Cells.Item(136, 8).Value = 23280.125
$ws.Cells.Item(136, 9).Value = 21968.4
$ws.Cells.Item(136, 11).Value = 65905.20000000001
$ws.Cells.Item(136, 13).Value = -63355.20000000001
